# GDV-Wochenplan.xlsx - add "Meilenstein" row (row 7) to Tabelle1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New milestone labels in row 7
$ws.Range("H7").Value = "Meilenstein 1"
$ws.Range("K7").Value = "Meilenstein 2"
$ws.Range("N7").Value = "Meilenstein 3"

# Widen column N (14) slightly to fit the new label
$ws.Columns(14).ColumnWidth = 14.57

# Move/restore the active selection to M9
$ws.Range("M9").Select() | Out-Null
